$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2,4,5,8,9,11,12,14,17,18,19,20,22,23)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = $false
}
